$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[William G.%Lindsley%NULL%1,                  Stephen B.%Martin%NULL%1,                  Robert E.%Thewlis%NULL%1,                  Khachatur%Sarkisian%NULL%1,                  Julian O.%Nwoko%NULL%1,                  Kenneth R.%Mead%NULL%1,                  John D.%Noti%NULL%1]'
$ws.Range("E3").Value = '[Michael S.%Bergman%xref no email%1,          Dennis J.%Viscusi%xref no email%1,          Brian K.%Heimbuch%xref no email%1,          Joseph D.%Wander%xref no email%1,          Anthony R.%Sambol%xref no email%1,          Ronald E.%Shaffer%xref no email%2]'
$ws.Range("E5").Value = '[Dennis J.%Viscusi%NULL%2,                  Michael S.%Bergman%NULL%1,                  Benjamin C.%Eimer%NULL%1,                  Ronald E.%Shaffer%rshaffer@cdc.gov%2]'
$ws.Range("E6").Value = '[Edward M.%Fisher%NULL%1,                  Jessica L.%Williams%NULL%1,                  Ronald E.%Shaffer%NULL%0,                  Tom%Jefferson%NULL%2,                  Tom%Jefferson%NULL%0]'
$ws.Range("E7").Value = '[ Brian K.%Heimbuch%null%2,                   Kimberly%Kinney%null%2,                   April E.%Lumley%null%2,                   Delbert A.%Harnish%null%1,                   Michael%Bergman%null%1,                   Joseph D.%Wander%null%2,               Brian K.%Heimbuch%null%2,               Kimberly%Kinney%null%2,               April E.%Lumley%null%2,               Delbert A.%Harnish%null%1,               Michael%Bergman%null%1,               Joseph D.%Wander%null%2]'
$ws.Range("E8").Value = '[Devin%Mills%NULL%1,                  Delbert A.%Harnish%dharnish@ara.com%1,                  Caryn%Lawrence%NULL%1,                  Megan%Sandoval-Powers%NULL%1,                  Brian K.%Heimbuch%NULL%1]'
$ws.Range("E9").Value = '[ Brian K.%Heimbuch%null%0,                   William H.%Wallace%null%1,                   Kimberly%Kinney%null%0,                   April E.%Lumley%null%0,                   Chang-Yu%Wu%null%1,                   Myung-Heui%Woo%null%1,                   Joseph D.%Wander%null%0,               Brian K.%Heimbuch%null%0,               William H.%Wallace%null%1,               Kimberly%Kinney%null%0,               April E.%Lumley%null%0,               Chang-Yu%Wu%null%1,               Myung-Heui%Woo%null%1,               Joseph D.%Wander%null%0]'
$ws.Range("E10").Value = '[Dennis J.%Viscusi%xref no email%0,          Michael S.%Bergman%xref no email%1,          Debra A.%Novak%xref no email%1,          Kimberly A.%Faulkner%xref no email%1,          Andrew%Palmiero%xref no email%1,          Jeffrey%Powell%xref no email%1,          Ronald E.%Shaffer%xref no email%0]'
$ws.Range("E11").Value = '[Diptanu%Paul%NULL%2,                  Ayush%Gupta%NULL%2,                  Anand Kumar%Maurya%NULL%4,                  Anand Kumar%Maurya%NULL%0,                  Amitava%Mukherjee%NULL%22,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0]'
$ws.Range("E12").Value = '[Diptanu%Paul%NULL%0,                  Ayush%Gupta%NULL%0,                  Anand Kumar%Maurya%NULL%0,                  Anand Kumar%Maurya%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0]'
$ws.Range("E13").Value = '[E.M.%Fisher%NULL%1,                  R.E.%Shaffer%NULL%1]'
$ws.Range("E14").Value = '[Evanly%Vo%xref no email%1,          Samy%Rengasamy%xref no email%2,          Ronald%Shaffer%xref no email%2]'
$ws.Range("E15").Value = '[Edward%Fisher%xref no email%1,          Samy%Rengasamy%xref no email%0,          Dennis%Viscusi%xref no email%1,          Evanly%Vo%xref no email%1,          Ronald%Shaffer%xref no email%0]'
$ws.Range("E16").Value = '[W.B.%Salter%NULL%1,                  K.%Kinney%NULL%1,                  W.H.%Wallace%NULL%1,                  A.E.%Lumley%NULL%1,                  B.K.%Heimbuch%NULL%1,                  J.D.%Wander%NULL%1]'
$ws.Range("E17").Value = '[Tzu-Hsien%Lin%NULL%1,                  Chih-Chieh%Chen%NULL%1,                  Sheng-Hsiu%Huang%NULL%1,                  Chung-Wen%Kuo%NULL%2,                  Chung-Wen%Kuo%NULL%0,                  Chane-Yu%Lai%NULL%1,                  Wen-Yinn%Lin%NULL%1,                  Amitava%Mukherjee%NULL%0,                  Amitava%Mukherjee%NULL%0]'
$ws.Range("E18").Value = '[T.‐H.%Lin%NULL%1,                  F.‐C.%Tang%NULL%1,                  P.‐C.%Hung%NULL%1,                  Z.‐C.%Hua%NULL%1,                  C.‐Y.%Lai%cylai511@gmail.com%1]'
